# Agile Product Backlog update - "Update up to sprint 22"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Agile Product Backlog")

# --- Index shift fix-ups caused by removing/renaming obsolete task names ---
$ws.Range("C33").Value = "Implement Confirm_required_contact_info_provided method  in User account web form"
$ws.Range("C46").Value = "Implement method Delete_All_Job_Problems in Jobs webform"
$ws.Range("C48").Value = "Implement Delete_All_Job_Problems in Jobs webform"
$ws.Range("C49").Value = "Implement DeleteAllJobChats in Jobs webform"
$ws.Range("C50").Value = "Implement ReadJob in Jobs webform"
$ws.Range("C54").Value = "Implement Auto_Insert_Job_Start_Date webform"
$ws.Range("C56").Value = "Implement Confirm_required_job_info_provided method  in Job web form"
$ws.Range("C57").Value = "Implement Register_User_For_Job method in Assign Users webform"

# --- Sprint 15/16 (problem system) group & task renames ---
$ws.Range("C59").Value = "Create problem system"
$ws.Range("C60").Value = "Implement CreateProblem in Assign Jobs Problem webform"
$ws.Range("C61").Value = "Implement UpdateProblem method in Jobs Problem webform"
$ws.Range("C62").Value = "Implement DeleteProblem method in Jobs Problem webform"
$ws.Range("C63").Value = "Implement reading and searching for problem system"
$ws.Range("C64").Value = "Implement ReadProblem method in Jobs Problem webform"
$ws.Range("C65").Value = "Implement ReadProblems method in Jobs Problem webform"
$ws.Range("C66").Value = "Implement SearchProblems method in jobs webform"

# --- Sprint 17 group header + fix typo on existing task ---
$ws.Range("C67").Value = "Validate required job data convert problem serverity to text"
$ws.Range("C68").Value = "Implement Confirm_required_problem_info_provided method in Job webform"
$ws.Range("C69").Value = "Implement  Convert_problem_severity_to_text method in Job Problem webform"

# Row 70 cleared (task removed)
$ws.Range("B70").Value = $null
$ws.Range("C70").Value = $null

# --- Sprint 18 group (chat CRUD) ---
$ws.Range("C71").Value = "Implement Attachment uploads for problem and perform Chat CRUD ops"
$ws.Range("B72").Value = 48
$ws.Range("C72").Value = "Implement Attach_zip_folder method in Job Problem webform"
$ws.Range("B73").Value = 49
$ws.Range("C73").Value = "Implement Attach_zip_file method in Job Problem webform"
$ws.Range("B74").Value = 50
$ws.Range("C74").Value = "Implement CreateChat method in Job Problem webform"

# --- Sprint 19 group (chat CRUD cont.) ---
$ws.Range("C75").Value = "Implement Chat CRUD ops"
$ws.Range("B76").Value = 51
$ws.Range("C76").Value = "Implement Publish_chat_time_and_date method in Job Problem webform"
$ws.Range("B77").Value = 52
$ws.Range("C77").Value = "Implement ReadChats method in Job Problem webform"
$ws.Range("B78").Value = 53
$ws.Range("C78").Value = "Implement UpdateChat method in Job Problem webform"

# --- Sprint 20 group (chat CRUD cont.) ---
$ws.Range("C79").Value = "Implement Chat CRUD ops and validate required data"
$ws.Range("B80").Value = 54
$ws.Range("C80").Value = "Implement DeleteChat method in Job Problem webform"
$ws.Range("B81").Value = 55
$ws.Range("C81").Value = "Implement Confirm_required_info_is_not_empty method in JobProblem webform"

# --- Sprint 21 group (testing) ---
$ws.Range("D83").Value = "Sprint 21"
$ws.Range("C83").Value = "Test application"
$ws.Range("B84").Value = 56
$ws.Range("C84").Value = "Input sample data into app."
$ws.Range("B85").Value = 57
$ws.Range("C85").Value = "Create release and create database and reset all primary keys to zero"

# --- Sprint 22 group (deployment) ---
$ws.Range("D87").Value = "Sprint 22"
$ws.Range("C87").Value = "Deployment"
$ws.Range("B88").Value = 58
$ws.Range("C88").Value = "Package app into zip"
$ws.Range("B89").Value = 59

# --- Placeholder sprint groups further down renumbered 23-26 ---
$ws.Range("D91").Value = "Sprint 23"
$ws.Range("D95").Value = "Sprint 24"
$ws.Range("D99").Value = "Sprint 25"
$ws.Range("D103").Value = "Sprint 26"

# --- Update frozen pane / active selection to reflect work further down the sheet ---
$ws.Application.ActiveWindow.ScrollRow = 76
$ws.Range("C89").Select()
